# Payroll Suite ScottishTaxWeek11 workbook edit:
# Rename the "Weekly" worksheets/labels to "Monthly" and renumber the
# "DO NOT TOUCH AUTOMATION EMP 107" marker down to "EMP 105".

$wb = $excel.ActiveWorkbook

$wsFirst   = $wb.Worksheets.Item("first")
$wsWeekly  = $wb.Worksheets.Item("GeneralTaxRateWeekly")
$wsProcess = $wb.Worksheets.Item("ProcessPayrollForWeeklyTax")
$wsReports = $wb.Worksheets.Item("TestReports")

# --- Rename the two "Weekly" sheets to their "Monthly" counterparts ---
$wsWeekly.Name = "GeneralTaxRateMonthly"
$wsProcess.Name = "ProcessPayrollForMonthlyTax"

# --- "first" sheet: the TC table references the (now renamed) sheet names ---
$wsFirst.Range("A3").Value = "GeneralTaxRateMonthly"
$wsFirst.Range("A4").Value = "ProcessPayrollForMonthlyTax"

# --- Update the "DO NOT TOUCH AUTOMATION EMP 107" marker -> "EMP 105" ---
$wsWeekly.Range("A2").Value = "DO NOT TOUCH AUTOMATION EMP 105"
$wsProcess.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 105"
$wsReports.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 105"

# --- Restore / update each sheet's selection, then activate sheets in order
#     so that the final active tab + per-sheet selections match the target.
$wsFirst.Activate() | Out-Null
$wsFirst.Range("F5").Select() | Out-Null

$wsWeekly.Activate() | Out-Null
$wsWeekly.Range("J2").Select() | Out-Null

$wsProcess.Activate() | Out-Null
$wsProcess.Range("E2:F2").Select() | Out-Null

$wsReports.Activate() | Out-Null
$wsReports.Range("K10").Select() | Out-Null
